# Update STATUS / COMMENT values across Sheet1 and Sheet2 per the QA
# Excel Compiler STATUS-tracking change:
#   ERROR -> ISSUE
#   OK    -> NO ISSUE
# plus a few COMMENT text simplifications.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("E3").Value = "ISSUE"
$ws1.Range("F3").Value = "Agree - typo"
$ws1.Range("E5").Value = "NO ISSUE"
$ws1.Range("F5").Value = "Good"
$ws1.Range("E6").Value = "ISSUE"

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("E2").Value = "NO ISSUE"
$ws2.Range("E3").Value = "NO ISSUE"
$ws2.Range("F3").Value = "Verified"
$ws2.Range("E4").Value = "NO ISSUE"
$ws2.Range("F4").Value = "Correct"
